# İş Takip Güncellemesi - 05.08.2025 14:39:48
# Fill in tracking dates / status values on the "Güncelleme" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Güncelleme")

# The source workbook stores these date-looking values as plain text
# (t="str"), not as real Excel dates. Force the target cells to text
# format first so COM doesn't auto-convert the "yyyy-mm-dd" strings into
# date serial numbers.
$dateCells = @("I3","K3","I4","K4","N4","O4","I6","K6","I8","K8","N8","O8","I10","K10","I23")
foreach ($addr in $dateCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 3
$ws.Range("I3").Value = "2024-11-05"
$ws.Range("K3").Value = "2024-12-11"
$ws.Range("L3").Value = "Yapıldı"
$ws.Range("M3").Value = "Yapıldı"

# Row 4
$ws.Range("I4").Value = "2024-11-05"
$ws.Range("K4").Value = "2024-12-11"
$ws.Range("L4").Value = "Yapıldı"
$ws.Range("M4").Value = "Yapıldı"
$ws.Range("N4").Value = "2025-05-13"
$ws.Range("O4").Value = "2025-05-13"

# Row 6
$ws.Range("I6").Value = "2024-11-07"
$ws.Range("K6").Value = "2024-12-13"
$ws.Range("L6").Value = "Yapıldı"
$ws.Range("M6").Value = "Yapıldı"

# Row 8
$ws.Range("I8").Value = "2024-11-07"
$ws.Range("K8").Value = "2024-12-11"
$ws.Range("L8").Value = "Yapıldı"
$ws.Range("M8").Value = "Yapıldı"
$ws.Range("N8").Value = "2025-06-03"
$ws.Range("O8").Value = "2025-06-03"

# Row 10
$ws.Range("I10").Value = "2024-11-07"
$ws.Range("K10").Value = "2024-12-13"
$ws.Range("L10").Value = "Yapıldı"
$ws.Range("M10").Value = "Yapıldı"

# Row 23
$ws.Range("I23").Value = "2024-11-11"
